$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet had a row-alignment bug: starting at row 13, the B/C "value"
# column content was one row ahead of the A "label" column. Fixing this
# shifts everything in rows 13-23 down by one row (to 14-24) and inserts
# new content that had been missing (the new row 13 label cell stays
# empty, row 10 B/C gets new "Objetivos" text, etc).
#
# Insert a blank row at row 13 - this shifts rows 13:23 down to 14:24,
# carrying row heights/styles along, which already lines up almost
# everything correctly.
$ws.Rows("13:13").Insert()

# New "Objetivos" (Objectives) long-form text replaces the professor
# name that was previously (incorrectly) sitting in B10/C10.
$ws.Range("B10").Value = "Complementar a formação multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, tópicos atuais e relevantes sobre fenômenos de transporte, termodinâmica, operações unitárias e reatore"
$ws.Range("C10").Value = "Complementar a formação multidisciplinar dos alunos de Engenharia abordando, com maior profundidade, tópicos atuais e relevantes sobre fenômenos de transporte, termodinâmica, operações unitárias e reatore"

# Row 13 ("Docentes responsáveis:" in A12, blank A13) now holds the
# professor name in B/C.
$ws.Range("B13").Value = "5816812 - João Paulo Alves Silva"
$ws.Range("C13").Value = "5816812 - João Paulo Alves Silva"

# Row 14 ("Programa resumido:") gets the real summary text instead of
# the leftover "Semestral" value.
$ws.Range("B14").Value = "Tópicos atuais e relevantes sobre fenômenos de transporte, termodinâmica, operações unitárias e reatores."
$ws.Range("C14").Value = "Tópicos atuais e relevantes sobre fenômenos de transporte, termodinâmica, operações unitárias e reatores."

# Row 15 ("Short syllabus:") has no B/C value - clear the stray date
# that the insert carried down.
$ws.Range("B15").Value = $null
$ws.Range("C15").Value = $null

# Row 16 ("Programa:") gets the real program description instead of
# the leftover date value.
$ws.Range("B16").Value = "O conteúdo desta disciplina será de acordo com o tópico a ser programado, devendo abordar assuntos complementares a formação de um profissional de Engenharia Química."
$ws.Range("C16").Value = "O conteúdo desta disciplina será de acordo com o tópico a ser programado, devendo abordar assuntos complementares a formação de um profissional de Engenharia Química."

# Row 18 ("Avaliação:") has no B/C value - clear the stray professor
# name that the insert carried down.
$ws.Range("B18").Value = $null
$ws.Range("C18").Value = $null

# Row 19 ("Método:") gets the real method description instead of the
# leftover professor name.
$ws.Range("B19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."
$ws.Range("C19").Value = "O desenvolvimento da disciplina será baseado em leituras, aula expositiva, discussão e resolução de estudos de caso e resolução de exercícios."

# Row 20 ("Critério:") gets "Provas e trabalhos." instead of the
# leftover method description.
$ws.Range("B20").Value = "Provas e trabalhos."
$ws.Range("C20").Value = "Provas e trabalhos."

# Row 21 ("Norma de recuperação:") gets the recovery-grade text instead
# of the leftover "Provas e trabalhos." value.
$ws.Range("B21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."
$ws.Range("C21").Value = "Prova única com nota maior ou igual a 5,0 (cinco)."

# Row 22 ("Bibliografia:") gets the real bibliography text instead of
# the leftover recovery-grade text.
$ws.Range("B22").Value = "Textos fornecidos pelo professor da disciplina`nArtigos extraídos de revistas especializadas de Engenharia Química."
$ws.Range("C22").Value = "Textos fornecidos pelo professor da disciplina`nArtigos extraídos de revistas especializadas de Engenharia Química."

# Rows 23 ("Requisitos:") / 24 (the prerequisite course text) already
# line up correctly after the insert - no further changes needed there.
